# Updates cryptos list values (price + volume%) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.254.31"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "3.480.09"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'610.15"
$ws.Range("E5").Value = "  +4.82%  "
$ws.Range("D6").Value = "'185.68"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.215"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").Value = "'0.651"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'52.99"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").Value = "'0.0000309"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "'9.50"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "4.034.00"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "'605.13"
$ws.Range("E15").Value = "  +6.22%  "
$ws.Range("D16").Value = "69.326.16"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.80"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.55"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "3.521.28"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'0.984"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").Value = "'17.14"
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("D23").Value = "'105.81"
$ws.Range("E23").Value = "  +10.84%  "
$ws.Range("D24").Value = "'4.63"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'5.05"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").Value = "'3.01"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "'9.79"
$ws.Range("E28").Value = "  +7.02%  "
$ws.Range("D29").Value = "'33.52"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").Value = "'6.94"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").Value = "'12.33"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'3.89"
$ws.Range("E33").Value = "  +15.42%  "
$ws.Range("D34").Value = "'63.10"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'3.21"
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'518.74"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("D40").Value = "3.578.30"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'36.68"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").Value = "0.0₃0775"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "'0.0460"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").Value = "'8.80"
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").Value = "'0.000244"
$ws.Range("E50").Value = "  -7.55%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "'1.36"
$ws.Range("E51").Value = "  -8.82%  "
